$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 56 - GenomeWeb article about Guardant Health / Eli Lilly Inluriyo CDx
$ws.Range("A56").Value = "https://www.genomeweb.com/cancer/guardant-health-blood-test-gets-fda-ok-cdx-eli-lilly-breast-cancer-drug-inluriyo"
$ws.Hyperlinks.Add($ws.Range("A56"), "https://www.genomeweb.com/cancer/guardant-health-blood-test-gets-fda-ok-cdx-eli-lilly-breast-cancer-drug-inluriyo")
$ws.Range("A56").Style = "Hyperlink"
$ws.Range("B56").Value = "CDx, companion diagnostic"
$ws.Range("C56").Value = "Guardant Health Blood Test Gets FDA OK as CDx for Eli Lilly Breast Cancer Drug Inluriyo"

# New row 57 - 360Dx mirror of the same article
$ws.Range("A57").Value = "https://www.360dx.com/cancer/guardant-health-blood-test-gets-fda-ok-cdx-eli-lilly-breast-cancer-drug-inluriyo"
$ws.Hyperlinks.Add($ws.Range("A57"), "https://www.360dx.com/cancer/guardant-health-blood-test-gets-fda-ok-cdx-eli-lilly-breast-cancer-drug-inluriyo")
$ws.Range("A57").Style = "Hyperlink"
$ws.Range("B57").Value = "CDx, companion diagnostic"
$ws.Range("C57").Value = "Guardant Health Blood Test Gets FDA OK as CDx for Eli Lilly Breast Cancer Drug Inluriyo"
